$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new row at position 21 (shifts old rows 21-94 down to 22-95)
$ws.Rows(21).Insert()

# Insert second new row at position 89 (after the first insert, old row 88 sits at 89;
# this pushes it - and everything below - down to 90, landing the new row at 89)
$ws.Rows(89).Insert()

# Populate the common (unchanged) columns for the two brand-new rows
foreach ($r in 21,89) {
    $ws.Cells.Item($r, 1).Value = 10
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value = "La Araucanía"
    $ws.Cells.Item($r, 5).Value = 9
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108007
    $ws.Cells.Item($r, 10).Value = "Coco"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = "Primera"
    $ws.Cells.Item($r, 17).Value = "`$/malla 20 unidades"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 20).Value = 20
}

# Set the per-row varying columns (D = Fecha, M = Volumen, N/O/P = precios, S = Precio $/Kg)
# for every data row 21..96 to match the final (shifted) dataset.
$ws.Cells.Item(21, 4).Value = "2022-10-25"
$ws.Cells.Item(21, 13).Value = 20
$ws.Cells.Item(21, 14).Value = 32000
$ws.Cells.Item(21, 15).Value = 32000
$ws.Cells.Item(21, 16).Value = 32000
$ws.Cells.Item(21, 19).Value = 1600
$ws.Cells.Item(22, 4).Value = "2022-08-03"
$ws.Cells.Item(22, 13).Value = 40
$ws.Cells.Item(22, 14).Value = 30000
$ws.Cells.Item(22, 15).Value = 30000
$ws.Cells.Item(22, 16).Value = 30000
$ws.Cells.Item(22, 19).Value = 1500
$ws.Cells.Item(23, 4).Value = "2022-06-03"
$ws.Cells.Item(23, 13).Value = 10
$ws.Cells.Item(23, 14).Value = 30000
$ws.Cells.Item(23, 15).Value = 30000
$ws.Cells.Item(23, 16).Value = 30000
$ws.Cells.Item(23, 19).Value = 1500
$ws.Cells.Item(24, 4).Value = "2022-06-15"
$ws.Cells.Item(24, 13).Value = 40
$ws.Cells.Item(24, 14).Value = 28000
$ws.Cells.Item(24, 15).Value = 28000
$ws.Cells.Item(24, 16).Value = 28000
$ws.Cells.Item(24, 19).Value = 1400
$ws.Cells.Item(25, 4).Value = "2022-04-20"
$ws.Cells.Item(25, 13).Value = 25
$ws.Cells.Item(25, 14).Value = 30000
$ws.Cells.Item(25, 15).Value = 30000
$ws.Cells.Item(25, 16).Value = 30000
$ws.Cells.Item(25, 19).Value = 1500
$ws.Cells.Item(26, 4).Value = "2022-06-28"
$ws.Cells.Item(26, 13).Value = 40
$ws.Cells.Item(26, 14).Value = 28000
$ws.Cells.Item(26, 15).Value = 28000
$ws.Cells.Item(26, 16).Value = 28000
$ws.Cells.Item(26, 19).Value = 1400
$ws.Cells.Item(27, 4).Value = "2022-05-18"
$ws.Cells.Item(27, 13).Value = 20
$ws.Cells.Item(27, 14).Value = 32000
$ws.Cells.Item(27, 15).Value = 32000
$ws.Cells.Item(27, 16).Value = 32000
$ws.Cells.Item(27, 19).Value = 1600
$ws.Cells.Item(28, 4).Value = "2021-02-08"
$ws.Cells.Item(28, 13).Value = 15
$ws.Cells.Item(28, 14).Value = 25000
$ws.Cells.Item(28, 15).Value = 25000
$ws.Cells.Item(28, 16).Value = 25000
$ws.Cells.Item(28, 19).Value = 1250
$ws.Cells.Item(29, 4).Value = "2020-12-10"
$ws.Cells.Item(29, 13).Value = 25
$ws.Cells.Item(29, 14).Value = 23000
$ws.Cells.Item(29, 15).Value = 23000
$ws.Cells.Item(29, 16).Value = 23000
$ws.Cells.Item(29, 19).Value = 1150
$ws.Cells.Item(30, 4).Value = "2022-06-08"
$ws.Cells.Item(30, 13).Value = 20
$ws.Cells.Item(30, 14).Value = 28000
$ws.Cells.Item(30, 15).Value = 28000
$ws.Cells.Item(30, 16).Value = 28000
$ws.Cells.Item(30, 19).Value = 1400
$ws.Cells.Item(31, 4).Value = "2022-07-13"
$ws.Cells.Item(31, 13).Value = 25
$ws.Cells.Item(31, 14).Value = 28000
$ws.Cells.Item(31, 15).Value = 28000
$ws.Cells.Item(31, 16).Value = 28000
$ws.Cells.Item(31, 19).Value = 1400
$ws.Cells.Item(32, 4).Value = "2021-08-10"
$ws.Cells.Item(32, 13).Value = 20
$ws.Cells.Item(32, 14).Value = 24000
$ws.Cells.Item(32, 15).Value = 24000
$ws.Cells.Item(32, 16).Value = 24000
$ws.Cells.Item(32, 19).Value = 1200
$ws.Cells.Item(33, 4).Value = "2022-08-05"
$ws.Cells.Item(33, 13).Value = 15
$ws.Cells.Item(33, 14).Value = 30000
$ws.Cells.Item(33, 15).Value = 30000
$ws.Cells.Item(33, 16).Value = 30000
$ws.Cells.Item(33, 19).Value = 1500
$ws.Cells.Item(34, 4).Value = "2021-07-13"
$ws.Cells.Item(34, 13).Value = 10
$ws.Cells.Item(34, 14).Value = 24000
$ws.Cells.Item(34, 15).Value = 24000
$ws.Cells.Item(34, 16).Value = 24000
$ws.Cells.Item(34, 19).Value = 1200
$ws.Cells.Item(35, 4).Value = "2021-07-23"
$ws.Cells.Item(35, 13).Value = 5
$ws.Cells.Item(35, 14).Value = 24000
$ws.Cells.Item(35, 15).Value = 24000
$ws.Cells.Item(35, 16).Value = 24000
$ws.Cells.Item(35, 19).Value = 1200
$ws.Cells.Item(36, 4).Value = "2021-07-15"
$ws.Cells.Item(36, 13).Value = 10
$ws.Cells.Item(36, 14).Value = 24000
$ws.Cells.Item(36, 15).Value = 24000
$ws.Cells.Item(36, 16).Value = 24000
$ws.Cells.Item(36, 19).Value = 1200
$ws.Cells.Item(37, 4).Value = "2021-09-22"
$ws.Cells.Item(37, 13).Value = 30
$ws.Cells.Item(37, 14).Value = 24000
$ws.Cells.Item(37, 15).Value = 24000
$ws.Cells.Item(37, 16).Value = 24000
$ws.Cells.Item(37, 19).Value = 1200
$ws.Cells.Item(38, 4).Value = "2021-09-29"
$ws.Cells.Item(38, 13).Value = 20
$ws.Cells.Item(38, 14).Value = 24000
$ws.Cells.Item(38, 15).Value = 24000
$ws.Cells.Item(38, 16).Value = 24000
$ws.Cells.Item(38, 19).Value = 1200
$ws.Cells.Item(39, 4).Value = "2021-08-18"
$ws.Cells.Item(39, 13).Value = 15
$ws.Cells.Item(39, 14).Value = 24000
$ws.Cells.Item(39, 15).Value = 24000
$ws.Cells.Item(39, 16).Value = 24000
$ws.Cells.Item(39, 19).Value = 1200
$ws.Cells.Item(40, 4).Value = "2021-12-02"
$ws.Cells.Item(40, 13).Value = 20
$ws.Cells.Item(40, 14).Value = 28000
$ws.Cells.Item(40, 15).Value = 28000
$ws.Cells.Item(40, 16).Value = 28000
$ws.Cells.Item(40, 19).Value = 1400
$ws.Cells.Item(41, 4).Value = "2021-08-20"
$ws.Cells.Item(41, 13).Value = 15
$ws.Cells.Item(41, 14).Value = 24000
$ws.Cells.Item(41, 15).Value = 24000
$ws.Cells.Item(41, 16).Value = 24000
$ws.Cells.Item(41, 19).Value = 1200
$ws.Cells.Item(42, 4).Value = "2021-10-05"
$ws.Cells.Item(42, 13).Value = 20
$ws.Cells.Item(42, 14).Value = 24000
$ws.Cells.Item(42, 15).Value = 24000
$ws.Cells.Item(42, 16).Value = 24000
$ws.Cells.Item(42, 19).Value = 1200
$ws.Cells.Item(43, 4).Value = "2021-02-22"
$ws.Cells.Item(43, 13).Value = 15
$ws.Cells.Item(43, 14).Value = 25000
$ws.Cells.Item(43, 15).Value = 25000
$ws.Cells.Item(43, 16).Value = 25000
$ws.Cells.Item(43, 19).Value = 1250
$ws.Cells.Item(44, 4).Value = "2022-04-12"
$ws.Cells.Item(44, 13).Value = 20
$ws.Cells.Item(44, 14).Value = 28000
$ws.Cells.Item(44, 15).Value = 28000
$ws.Cells.Item(44, 16).Value = 28000
$ws.Cells.Item(44, 19).Value = 1400
$ws.Cells.Item(45, 4).Value = "2022-05-23"
$ws.Cells.Item(45, 13).Value = 12
$ws.Cells.Item(45, 14).Value = 30000
$ws.Cells.Item(45, 15).Value = 30000
$ws.Cells.Item(45, 16).Value = 30000
$ws.Cells.Item(45, 19).Value = 1500
$ws.Cells.Item(46, 4).Value = "2022-04-13"
$ws.Cells.Item(46, 13).Value = 25
$ws.Cells.Item(46, 14).Value = 28000
$ws.Cells.Item(46, 15).Value = 28000
$ws.Cells.Item(46, 16).Value = 28000
$ws.Cells.Item(46, 19).Value = 1400
$ws.Cells.Item(47, 4).Value = "2021-08-13"
$ws.Cells.Item(47, 13).Value = 20
$ws.Cells.Item(47, 14).Value = 24000
$ws.Cells.Item(47, 15).Value = 24000
$ws.Cells.Item(47, 16).Value = 24000
$ws.Cells.Item(47, 19).Value = 1200
$ws.Cells.Item(48, 4).Value = "2022-07-08"
$ws.Cells.Item(48, 13).Value = 30
$ws.Cells.Item(48, 14).Value = 28000
$ws.Cells.Item(48, 15).Value = 28000
$ws.Cells.Item(48, 16).Value = 28000
$ws.Cells.Item(48, 19).Value = 1400
$ws.Cells.Item(49, 4).Value = "2022-08-04"
$ws.Cells.Item(49, 13).Value = 30
$ws.Cells.Item(49, 14).Value = 30000
$ws.Cells.Item(49, 15).Value = 30000
$ws.Cells.Item(49, 16).Value = 30000
$ws.Cells.Item(49, 19).Value = 1500
$ws.Cells.Item(50, 4).Value = "2022-06-07"
$ws.Cells.Item(50, 13).Value = 25
$ws.Cells.Item(50, 14).Value = 28000
$ws.Cells.Item(50, 15).Value = 28000
$ws.Cells.Item(50, 16).Value = 28000
$ws.Cells.Item(50, 19).Value = 1400
$ws.Cells.Item(51, 4).Value = "2022-07-14"
$ws.Cells.Item(51, 13).Value = 25
$ws.Cells.Item(51, 14).Value = 28000
$ws.Cells.Item(51, 15).Value = 28000
$ws.Cells.Item(51, 16).Value = 28000
$ws.Cells.Item(51, 19).Value = 1400
$ws.Cells.Item(52, 4).Value = "2021-06-09"
$ws.Cells.Item(52, 13).Value = 15
$ws.Cells.Item(52, 14).Value = 24000
$ws.Cells.Item(52, 15).Value = 24000
$ws.Cells.Item(52, 16).Value = 24000
$ws.Cells.Item(52, 19).Value = 1200
$ws.Cells.Item(53, 4).Value = "2021-08-06"
$ws.Cells.Item(53, 13).Value = 15
$ws.Cells.Item(53, 14).Value = 25000
$ws.Cells.Item(53, 15).Value = 25000
$ws.Cells.Item(53, 16).Value = 25000
$ws.Cells.Item(53, 19).Value = 1250
$ws.Cells.Item(54, 4).Value = "2022-07-04"
$ws.Cells.Item(54, 13).Value = 70
$ws.Cells.Item(54, 14).Value = 28000
$ws.Cells.Item(54, 15).Value = 30000
$ws.Cells.Item(54, 16).Value = 29143
$ws.Cells.Item(54, 19).Value = 1457
$ws.Cells.Item(55, 4).Value = "2021-07-14"
$ws.Cells.Item(55, 13).Value = 10
$ws.Cells.Item(55, 14).Value = 24000
$ws.Cells.Item(55, 15).Value = 24000
$ws.Cells.Item(55, 16).Value = 24000
$ws.Cells.Item(55, 19).Value = 1200
$ws.Cells.Item(56, 4).Value = "2022-05-20"
$ws.Cells.Item(56, 13).Value = 50
$ws.Cells.Item(56, 14).Value = 30000
$ws.Cells.Item(56, 15).Value = 32000
$ws.Cells.Item(56, 16).Value = 31200
$ws.Cells.Item(56, 19).Value = 1560
$ws.Cells.Item(57, 4).Value = "2021-08-11"
$ws.Cells.Item(57, 13).Value = 40
$ws.Cells.Item(57, 14).Value = 25000
$ws.Cells.Item(57, 15).Value = 25000
$ws.Cells.Item(57, 16).Value = 25000
$ws.Cells.Item(57, 19).Value = 1250
$ws.Cells.Item(58, 4).Value = "2021-06-16"
$ws.Cells.Item(58, 13).Value = 30
$ws.Cells.Item(58, 14).Value = 24000
$ws.Cells.Item(58, 15).Value = 24000
$ws.Cells.Item(58, 16).Value = 24000
$ws.Cells.Item(58, 19).Value = 1200
$ws.Cells.Item(59, 4).Value = "2022-05-09"
$ws.Cells.Item(59, 13).Value = 35
$ws.Cells.Item(59, 14).Value = 30000
$ws.Cells.Item(59, 15).Value = 30000
$ws.Cells.Item(59, 16).Value = 30000
$ws.Cells.Item(59, 19).Value = 1500
$ws.Cells.Item(60, 4).Value = "2021-09-30"
$ws.Cells.Item(60, 13).Value = 40
$ws.Cells.Item(60, 14).Value = 24000
$ws.Cells.Item(60, 15).Value = 24000
$ws.Cells.Item(60, 16).Value = 24000
$ws.Cells.Item(60, 19).Value = 1200
$ws.Cells.Item(61, 4).Value = "2022-03-25"
$ws.Cells.Item(61, 13).Value = 10
$ws.Cells.Item(61, 14).Value = 28000
$ws.Cells.Item(61, 15).Value = 28000
$ws.Cells.Item(61, 16).Value = 28000
$ws.Cells.Item(61, 19).Value = 1400
$ws.Cells.Item(62, 4).Value = "2020-12-29"
$ws.Cells.Item(62, 13).Value = 20
$ws.Cells.Item(62, 14).Value = 20000
$ws.Cells.Item(62, 15).Value = 20000
$ws.Cells.Item(62, 16).Value = 20000
$ws.Cells.Item(62, 19).Value = 1000
$ws.Cells.Item(63, 4).Value = "2022-07-07"
$ws.Cells.Item(63, 13).Value = 100
$ws.Cells.Item(63, 14).Value = 28000
$ws.Cells.Item(63, 15).Value = 30000
$ws.Cells.Item(63, 16).Value = 28800
$ws.Cells.Item(63, 19).Value = 1440
$ws.Cells.Item(64, 4).Value = "2022-04-06"
$ws.Cells.Item(64, 13).Value = 20
$ws.Cells.Item(64, 14).Value = 28000
$ws.Cells.Item(64, 15).Value = 28000
$ws.Cells.Item(64, 16).Value = 28000
$ws.Cells.Item(64, 19).Value = 1400
$ws.Cells.Item(65, 4).Value = "2021-09-13"
$ws.Cells.Item(65, 13).Value = 25
$ws.Cells.Item(65, 14).Value = 25000
$ws.Cells.Item(65, 15).Value = 25000
$ws.Cells.Item(65, 16).Value = 25000
$ws.Cells.Item(65, 19).Value = 1250
$ws.Cells.Item(66, 4).Value = "2021-08-27"
$ws.Cells.Item(66, 13).Value = 100
$ws.Cells.Item(66, 14).Value = 24000
$ws.Cells.Item(66, 15).Value = 24000
$ws.Cells.Item(66, 16).Value = 24000
$ws.Cells.Item(66, 19).Value = 1200
$ws.Cells.Item(67, 4).Value = "2021-06-02"
$ws.Cells.Item(67, 13).Value = 30
$ws.Cells.Item(67, 14).Value = 24000
$ws.Cells.Item(67, 15).Value = 24000
$ws.Cells.Item(67, 16).Value = 24000
$ws.Cells.Item(67, 19).Value = 1200
$ws.Cells.Item(68, 4).Value = "2021-07-12"
$ws.Cells.Item(68, 13).Value = 20
$ws.Cells.Item(68, 14).Value = 24000
$ws.Cells.Item(68, 15).Value = 24000
$ws.Cells.Item(68, 16).Value = 24000
$ws.Cells.Item(68, 19).Value = 1200
$ws.Cells.Item(69, 4).Value = "2021-08-04"
$ws.Cells.Item(69, 13).Value = 20
$ws.Cells.Item(69, 14).Value = 25000
$ws.Cells.Item(69, 15).Value = 25000
$ws.Cells.Item(69, 16).Value = 25000
$ws.Cells.Item(69, 19).Value = 1250
$ws.Cells.Item(70, 4).Value = "2022-07-06"
$ws.Cells.Item(70, 13).Value = 40
$ws.Cells.Item(70, 14).Value = 28000
$ws.Cells.Item(70, 15).Value = 28000
$ws.Cells.Item(70, 16).Value = 28000
$ws.Cells.Item(70, 19).Value = 1400
$ws.Cells.Item(71, 4).Value = "2021-02-24"
$ws.Cells.Item(71, 13).Value = 15
$ws.Cells.Item(71, 14).Value = 25000
$ws.Cells.Item(71, 15).Value = 25000
$ws.Cells.Item(71, 16).Value = 25000
$ws.Cells.Item(71, 19).Value = 1250
$ws.Cells.Item(72, 4).Value = "2022-08-10"
$ws.Cells.Item(72, 13).Value = 35
$ws.Cells.Item(72, 14).Value = 30000
$ws.Cells.Item(72, 15).Value = 30000
$ws.Cells.Item(72, 16).Value = 30000
$ws.Cells.Item(72, 19).Value = 1500
$ws.Cells.Item(73, 4).Value = "2021-09-15"
$ws.Cells.Item(73, 13).Value = 25
$ws.Cells.Item(73, 14).Value = 25000
$ws.Cells.Item(73, 15).Value = 25000
$ws.Cells.Item(73, 16).Value = 25000
$ws.Cells.Item(73, 19).Value = 1250
$ws.Cells.Item(74, 4).Value = "2021-06-30"
$ws.Cells.Item(74, 13).Value = 15
$ws.Cells.Item(74, 14).Value = 20000
$ws.Cells.Item(74, 15).Value = 20000
$ws.Cells.Item(74, 16).Value = 20000
$ws.Cells.Item(74, 19).Value = 1000
$ws.Cells.Item(75, 4).Value = "2021-07-05"
$ws.Cells.Item(75, 13).Value = 15
$ws.Cells.Item(75, 14).Value = 20000
$ws.Cells.Item(75, 15).Value = 20000
$ws.Cells.Item(75, 16).Value = 20000
$ws.Cells.Item(75, 19).Value = 1000
$ws.Cells.Item(76, 4).Value = "2022-05-25"
$ws.Cells.Item(76, 13).Value = 25
$ws.Cells.Item(76, 14).Value = 30000
$ws.Cells.Item(76, 15).Value = 30000
$ws.Cells.Item(76, 16).Value = 30000
$ws.Cells.Item(76, 19).Value = 1500
$ws.Cells.Item(77, 4).Value = "2022-07-05"
$ws.Cells.Item(77, 13).Value = 20
$ws.Cells.Item(77, 14).Value = 28000
$ws.Cells.Item(77, 15).Value = 28000
$ws.Cells.Item(77, 16).Value = 28000
$ws.Cells.Item(77, 19).Value = 1400
$ws.Cells.Item(78, 4).Value = "2022-07-27"
$ws.Cells.Item(78, 13).Value = 20
$ws.Cells.Item(78, 14).Value = 30000
$ws.Cells.Item(78, 15).Value = 30000
$ws.Cells.Item(78, 16).Value = 30000
$ws.Cells.Item(78, 19).Value = 1500
$ws.Cells.Item(79, 4).Value = "2021-05-18"
$ws.Cells.Item(79, 13).Value = 20
$ws.Cells.Item(79, 14).Value = 25000
$ws.Cells.Item(79, 15).Value = 25000
$ws.Cells.Item(79, 16).Value = 25000
$ws.Cells.Item(79, 19).Value = 1250
$ws.Cells.Item(80, 4).Value = "2021-10-19"
$ws.Cells.Item(80, 13).Value = 40
$ws.Cells.Item(80, 14).Value = 20000
$ws.Cells.Item(80, 15).Value = 20000
$ws.Cells.Item(80, 16).Value = 20000
$ws.Cells.Item(80, 19).Value = 1000
$ws.Cells.Item(81, 4).Value = "2021-01-25"
$ws.Cells.Item(81, 13).Value = 30
$ws.Cells.Item(81, 14).Value = 25000
$ws.Cells.Item(81, 15).Value = 25000
$ws.Cells.Item(81, 16).Value = 25000
$ws.Cells.Item(81, 19).Value = 1250
$ws.Cells.Item(82, 4).Value = "2021-07-19"
$ws.Cells.Item(82, 13).Value = 12
$ws.Cells.Item(82, 14).Value = 24000
$ws.Cells.Item(82, 15).Value = 24000
$ws.Cells.Item(82, 16).Value = 24000
$ws.Cells.Item(82, 19).Value = 1200
$ws.Cells.Item(83, 4).Value = "2021-01-18"
$ws.Cells.Item(83, 13).Value = 15
$ws.Cells.Item(83, 14).Value = 25000
$ws.Cells.Item(83, 15).Value = 25000
$ws.Cells.Item(83, 16).Value = 25000
$ws.Cells.Item(83, 19).Value = 1250
$ws.Cells.Item(84, 4).Value = "2022-06-17"
$ws.Cells.Item(84, 13).Value = 20
$ws.Cells.Item(84, 14).Value = 28000
$ws.Cells.Item(84, 15).Value = 28000
$ws.Cells.Item(84, 16).Value = 28000
$ws.Cells.Item(84, 19).Value = 1400
$ws.Cells.Item(85, 4).Value = "2021-07-21"
$ws.Cells.Item(85, 13).Value = 15
$ws.Cells.Item(85, 14).Value = 25000
$ws.Cells.Item(85, 15).Value = 25000
$ws.Cells.Item(85, 16).Value = 25000
$ws.Cells.Item(85, 19).Value = 1250
$ws.Cells.Item(86, 4).Value = "2021-01-26"
$ws.Cells.Item(86, 13).Value = 15
$ws.Cells.Item(86, 14).Value = 25000
$ws.Cells.Item(86, 15).Value = 25000
$ws.Cells.Item(86, 16).Value = 25000
$ws.Cells.Item(86, 19).Value = 1250
$ws.Cells.Item(87, 4).Value = "2022-07-11"
$ws.Cells.Item(87, 13).Value = 20
$ws.Cells.Item(87, 14).Value = 25000
$ws.Cells.Item(87, 15).Value = 25000
$ws.Cells.Item(87, 16).Value = 25000
$ws.Cells.Item(87, 19).Value = 1250
$ws.Cells.Item(88, 4).Value = "2021-09-23"
$ws.Cells.Item(88, 13).Value = 10
$ws.Cells.Item(88, 14).Value = 24000
$ws.Cells.Item(88, 15).Value = 24000
$ws.Cells.Item(88, 16).Value = 24000
$ws.Cells.Item(88, 19).Value = 1200
$ws.Cells.Item(89, 4).Value = "2022-10-24"
$ws.Cells.Item(89, 13).Value = 30
$ws.Cells.Item(89, 14).Value = 32000
$ws.Cells.Item(89, 15).Value = 32000
$ws.Cells.Item(89, 16).Value = 32000
$ws.Cells.Item(89, 19).Value = 1600
$ws.Cells.Item(90, 4).Value = "2021-08-26"
$ws.Cells.Item(90, 13).Value = 20
$ws.Cells.Item(90, 14).Value = 24000
$ws.Cells.Item(90, 15).Value = 24000
$ws.Cells.Item(90, 16).Value = 24000
$ws.Cells.Item(90, 19).Value = 1200
$ws.Cells.Item(91, 4).Value = "2021-08-25"
$ws.Cells.Item(91, 13).Value = 10
$ws.Cells.Item(91, 14).Value = 24000
$ws.Cells.Item(91, 15).Value = 24000
$ws.Cells.Item(91, 16).Value = 24000
$ws.Cells.Item(91, 19).Value = 1200
$ws.Cells.Item(92, 4).Value = "2022-08-02"
$ws.Cells.Item(92, 13).Value = 20
$ws.Cells.Item(92, 14).Value = 30000
$ws.Cells.Item(92, 15).Value = 30000
$ws.Cells.Item(92, 16).Value = 30000
$ws.Cells.Item(92, 19).Value = 1500
$ws.Cells.Item(93, 4).Value = "2022-06-24"
$ws.Cells.Item(93, 13).Value = 20
$ws.Cells.Item(93, 14).Value = 28000
$ws.Cells.Item(93, 15).Value = 28000
$ws.Cells.Item(93, 16).Value = 28000
$ws.Cells.Item(93, 19).Value = 1400
$ws.Cells.Item(94, 4).Value = "2022-05-26"
$ws.Cells.Item(94, 13).Value = 25
$ws.Cells.Item(94, 14).Value = 30000
$ws.Cells.Item(94, 15).Value = 30000
$ws.Cells.Item(94, 16).Value = 30000
$ws.Cells.Item(94, 19).Value = 1500
$ws.Cells.Item(95, 4).Value = "2022-06-10"
$ws.Cells.Item(95, 13).Value = 15
$ws.Cells.Item(95, 14).Value = 20000
$ws.Cells.Item(95, 15).Value = 20000
$ws.Cells.Item(95, 16).Value = 20000
$ws.Cells.Item(95, 19).Value = 1000
$ws.Cells.Item(96, 4).Value = "2021-10-20"
$ws.Cells.Item(96, 13).Value = 40
$ws.Cells.Item(96, 14).Value = 24000
$ws.Cells.Item(96, 15).Value = 24000
$ws.Cells.Item(96, 16).Value = 24000
$ws.Cells.Item(96, 19).Value = 1200
